$d = $word.ActiveDocument

$p2 = $d.Paragraphs(2)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
          <w:i/>
        </w:rPr>
        <w:t>I think the cellar is on the right side</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t xml:space="preserve">, you think to yourself. The entire hallway is now empty so you run to the right side. </w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p2.Range.InsertXML($xml2)

$p3 = $d.Paragraphs(3)
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t xml:space="preserve">“Is it?” Narrator grins. </w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p3.Range.InsertXML($xml3)

$p4 = $d.Paragraphs(4)
$xml4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
          <w:i/>
        </w:rPr>
        <w:t>It should be one of these rooms.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t xml:space="preserve"> You spot a door with stairs on the other side leading what appears downward. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
          <w:i/>
        </w:rPr>
        <w:t>Aha! That’s got to be it!</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p4.Range.InsertXML($xml4)

$p5 = $d.Paragraphs(5)
$xml5 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t xml:space="preserve">You run down the stairs to see a door and open it. Jail cells greet you with welcoming arms. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
          <w:i/>
        </w:rPr>
        <w:t>This isn’t the cellar….</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t xml:space="preserve"> you realize.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p5.Range.InsertXML($xml5)

$p6 = $d.Paragraphs(6)
$xml6 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t>“Capture him!” You feel rough hands grabbing at your torso. Unable to shake them off you fall down face first. Your arms are pulled behind your back and locked together using a handcuff.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p6.Range.InsertXML($xml6)

$p7 = $d.Paragraphs(7)
$xml7 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
          <w:i/>
        </w:rPr>
        <w:t xml:space="preserve">Curses… I got caught. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t>The guards lock you inside one of the cells and take the ring from you.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p7.Range.InsertXML($xml7)

$p8 = $d.Paragraphs(8)
$xml8 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t>“You si</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t>t tight there young lad, I’ll get the sir.” The guards exit the room and you sit down with a sigh. First mission failed.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p8.Range.InsertXML($xml8)

$p9 = $d.Paragraphs(9)
$xml9 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t>“Nice going there,” Narrator snickers.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p9.Range.InsertXML($xml9)

$p10 = $d.Paragraphs(10)
$xml10 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t>“Shut up.”</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p10.Range.InsertXML($xml10)

$p11 = $d.Paragraphs(11)
$xml11 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t>“</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t>Pssst</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t xml:space="preserve">,” a voice whispers out from above you. You look up to see one of the thieves, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t>Rolin</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t xml:space="preserve"> poking his head out from the ceiling and gesturing you to follow him. You get up and jump up into the hole to follow him through the pathway he dug. Eventually </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t>the</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t xml:space="preserve"> two of you exit the tunnel. You find yourself in an alleyway. </w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p11.Range.InsertXML($xml11)

$p12 = $d.Paragraphs(12)
$xml12 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t>“So did you acquire</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t xml:space="preserve"> the ring?” </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t>Rolin</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t xml:space="preserve"> asks.</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p12.Range.InsertXML($xml12)

$p13 = $d.Paragraphs(13)
$xml13 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00A61CD2" w:rsidRDefault="00A61CD2" w:rsidP="00A61CD2">
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="3D3C40"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="3D3C40"/>
        </w:rPr>
        <w:t>--Answer truthfully</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$p13.Range.InsertXML($xml13)

$p14 = $d.Paragraphs(14)
$p15 = $d.Paragraphs(15)
$rng1415 = $d.Range($p14.Range.Start, $p15.Range.End)
$xml1415 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="3D3C40"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
          <w:b/>
          <w:bCs/>
          <w:color w:val="3D3C40"/>
        </w:rPr>
        <w:t>--Lie and pretend you have it</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Bookerly" w:hAnsi="Bookerly" w:cs="Bookerly"/>
        </w:rPr>
      </w:pPr>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$rng1415.InsertXML($xml1415)
